$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Locales"
$ws.Range("A2").Value = "Edificio sólo con una vivienda familiar"
$ws.Range("A3").Value = "Edificio sólo con varias viviendas familiares"
$ws.Range("A4").Value = "Edificios principalmente con vivienda colectiva: hotel, albergue, pensión"
$ws.Range("A5").Value = "Alojamientos"
$ws.Range("A6").Value = "Edificios principalmente con vivienda colectiva: convento, cuartel, prisión"
$ws.Range("A7").Value = "Edificios principalmente con vivienda colectiva: hospitales, instituciones para discapacitados,¿"
$ws.Range("A8").Value = "Edificios principalmente con locales compartidos con alguna vivienda"
$ws.Range("A9").Value = "Edificio principalmente con viviendas familiares compartido con locales"
$ws.Range("A10").Value = "Edificios principalmente con vivienda colectiva: instituciones de enseñanza,¿"
